# Regenerate the "K" column (column G) values for each save-data row.
# K replaces the old "Strike#" derived value; here we write the
# recalculated counts directly into column G for rows 2-47.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 0
    3  = 0
    4  = 0
    5  = 3
    6  = 1
    7  = 1
    8  = 1
    9  = 1
    10 = 1
    11 = 1
    12 = 2
    13 = 1
    14 = 1
    15 = 1
    16 = 0
    17 = 0
    18 = 1
    19 = 0
    20 = 1
    21 = 0
    22 = 0
    23 = 1
    24 = 1
    25 = 2
    26 = 1
    27 = 1
    28 = 0
    29 = 0
    30 = 2
    31 = 0
    32 = 0
    33 = 2
    34 = 0
    35 = 0
    36 = 0
    37 = 0
    38 = 0
    39 = 0
    40 = 1
    41 = 1
    42 = 0
    43 = 1
    44 = 0
    45 = 0
    46 = 0
    47 = 0
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
